$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "261.52"
Set-TextValue 2 5 "1.05%"
Set-TextValue 3 4 "27.22"
Set-TextValue 3 5 "1.46%"
Set-TextValue 4 4 "4.698"
Set-TextValue 4 5 "0.55%"
Set-TextValue 5 5 "3.24%"
Set-TextValue 6 4 "6.718"
Set-TextValue 6 5 "0.88%"
Set-TextValue 7 4 "0.8501"
Set-TextValue 7 5 "-0.91%"
Set-TextValue 8 4 "0.9157"
Set-TextValue 8 5 "-0.70%"
Set-TextValue 9 4 "0.1409"
Set-TextValue 9 5 "1.23%"
Set-TextValue 10 4 "0.04659"
Set-TextValue 10 5 "-5.72%"
Set-TextValue 11 4 "0.07090"
Set-TextValue 11 5 "1.16%"
Set-TextValue 12 4 "0.03131"
Set-TextValue 12 5 "2.85%"
Set-TextValue 13 4 "0.09040"
Set-TextValue 13 5 "-1.05%"
Set-TextValue 14 4 "0.001525"
Set-TextValue 14 5 "-0.88%"
Set-TextValue 15 4 "0.0006165"
Set-TextValue 15 5 "1.97%"
Set-TextValue 16 4 "0.006034"
Set-TextValue 16 5 "-1.11%"
Set-TextValue 17 4 "3.459"
Set-TextValue 17 5 "0.07%"
Set-TextValue 18 5 "0.62%"
Set-TextValue 19 4 "2.194"
Set-TextValue 19 5 "1.20%"
Set-TextValue 20 4 "0.3080"
Set-TextValue 20 5 "-0.95%"
Set-TextValue 21 4 "0.1308"
Set-TextValue 21 5 "1.56%"
Set-TextValue 22 4 "4.114"
Set-TextValue 22 5 "-0.72%"
Set-TextValue 23 4 "0.04227"
Set-TextValue 24 4 "0.001217"
Set-TextValue 24 5 "0.13%"
Set-TextValue 25 5 "-5.81%"
Set-TextValue 27 4 "0.0001600"
Set-TextValue 27 5 "-6.51%"
Set-TextValue 40 4 "0.03972"
Set-TextValue 40 5 "3.46%"
Set-TextValue 41 5 "-0.29%"
Set-TextValue 42 4 "0.004119"
Set-TextValue 42 5 "8.20%"
Set-TextValue 43 5 "-9.69%"
Set-TextValue 44 4 "0.01383"
Set-TextValue 44 5 "-8.29%"
Set-TextValue 45 4 "0.00005135"
Set-TextValue 45 5 "0.27%"
Set-TextValue 46 5 "0.09%"
Set-TextValue 48 4 "0.1667"
Set-TextValue 48 5 "10.88%"
Set-TextValue 49 4 "0.00002100"
Set-TextValue 49 5 "0.09%"
Set-TextValue 50 4 "0.0002000"
Set-TextValue 50 5 "0.09%"

Write-Host "Updated 63 cells"